# Updated cryptos list with GitHub Actions
# Applies the price/volume changes described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would be auto-parsed as a number by Excel need to be
# forced to Text format first (so they keep the exact literal formatting,
# e.g. "324.40" instead of becoming the number 324.4), then the number
# format is reset back to General / style back to Normal so no stray
# formatting is left behind on the cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "28.883.93"
$ws.Range("E2").Value = "  -0.29%  "
$ws.Range("D3").Value = "1.916.33"
$ws.Range("E3").Value = "  +0.71%  "
$ws.Range("E4").Value = "  +0.01%  "
Set-TextValue $ws.Range("D5") "324.40"
$ws.Range("E5").Value = "  -0.04%  "
Set-TextValue $ws.Range("D6") "1.003"
$ws.Range("E6").Value = "  -0.06%  "
Set-TextValue $ws.Range("D7") "0.4561"
$ws.Range("E7").Value = "  -0.74%  "
Set-TextValue $ws.Range("D8") "0.3797"
$ws.Range("E8").Value = "  -0.60%  "
Set-TextValue $ws.Range("D9") "0.07737"
$ws.Range("E9").Value = "  +0.37%  "
Set-TextValue $ws.Range("D10") "0.9743"
$ws.Range("E10").Value = "  -0.64%  "
Set-TextValue $ws.Range("D11") "22.25"
$ws.Range("E11").Value = "  +1.00%  "
$ws.Range("D12").Value = "1.949.72"
$ws.Range("E12").Value = "  -2.63%  "
Set-TextValue $ws.Range("D13") "5.688"
$ws.Range("E13").Value = "  +0.33%  "
Set-TextValue $ws.Range("D14") "6.965"
$ws.Range("E14").Value = "  -0.09%  "
Set-TextValue $ws.Range("D15") "0.06988"
$ws.Range("E15").Value = "  -0.81%  "
Set-TextValue $ws.Range("D16") "1.006"
$ws.Range("E16").Value = "  +0.11%  "
$ws.Range("E17").Value = "  +0.50%  "
Set-TextValue $ws.Range("D18") "0.000009460"
$ws.Range("E18").Value = "  -0.82%  "
$ws.Range("E19").Value = "  -0.76%  "
$ws.Range("E20").Value = "  -0.13%  "
$ws.Range("D21").Value = "28.908.72"
$ws.Range("E21").Value = "  -0.31%  "
Set-TextValue $ws.Range("D22") "5.329"
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("E23").Value = "  +1.46%  "
$ws.Range("D24").Value = "2.133.98"
$ws.Range("E24").Value = "  -3.22%  "
Set-TextValue $ws.Range("D25") "2.055"
$ws.Range("E25").Value = "  -2.17%  "
Set-TextValue $ws.Range("D26") "157.57"
$ws.Range("E26").Value = "  +0.13%  "
Set-TextValue $ws.Range("D27") "19.03"
$ws.Range("E27").Value = "  -0.42%  "
Set-TextValue $ws.Range("D28") "5.602"
$ws.Range("E28").Value = "  +0.18%  "
Set-TextValue $ws.Range("D29") "117.86"
$ws.Range("E29").Value = "  +0.26%  "
$ws.Range("E30").Value = "  -0.86%  "
Set-TextValue $ws.Range("D31") "0.09295"
$ws.Range("E31").Value = "  +0.14%  "
Set-TextValue $ws.Range("D32") "0.8661"
$ws.Range("E32").Value = "  +0.59%  "
$ws.Range("E33").Value = "  +0.19%  "
Set-TextValue $ws.Range("D34") "1.239"
$ws.Range("E34").Value = "  -0.93%  "
Set-TextValue $ws.Range("D35") "3.013"
$ws.Range("E35").Value = "  -0.12%  "
Set-TextValue $ws.Range("D36") "0.05680"
$ws.Range("E36").Value = "  -0.19%  "
$ws.Range("E37").Value = "  +0.07%  "
Set-TextValue $ws.Range("D38") "1.003"
$ws.Range("E38").Value = "  +0.03%  "
Set-TextValue $ws.Range("D39") "0.02036"
$ws.Range("E39").Value = "  +0.06%  "
Set-TextValue $ws.Range("D40") "3.066"
$ws.Range("E40").Value = "  +11.65%  "
Set-TextValue $ws.Range("D41") "7.452"
$ws.Range("E41").Value = "  -0.32%  "
$ws.Range("E42").Value = "  -0.58%  "
$ws.Range("E43").Value = "  -0.10%  "
Set-TextValue $ws.Range("D44") "9.303"
$ws.Range("E44").Value = "  +0.36%  "
Set-TextValue $ws.Range("D45") "0.000002818"
$ws.Range("E45").Value = "  +16.50%  "
Set-TextValue $ws.Range("D46") "2.154"
$ws.Range("E46").Value = "  +3.41%  "
$ws.Range("E47").Value = "  -0.90%  "
Set-TextValue $ws.Range("D48") "0.06929"
$ws.Range("E48").Value = "  +1.58%  "
$ws.Range("E49").Value = "  -1.59%  "
$ws.Range("E50").Value = "  -0.74%  "
$ws.Range("E51").Value = "  -1.04%  "
